$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 372
$ws.Range("F4").Value = 426
$ws.Range("F5").Value = 1155
$ws.Range("F8").Value = 1028
$ws.Range("F9").Value = 1635
$ws.Range("F10").Value = 6130
$ws.Range("G10").Value = 68
$ws.Range("F12").Value = 1770
$ws.Range("F13").Value = 459
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = 6060
$ws.Range("F16").Value = 6060
$ws.Range("F20").Value = 99
$ws.Range("F21").Value = 1668
$ws.Range("F25").Value = 1448
$ws.Range("F26").Value = 733
$ws.Range("F27").Value = 270
$ws.Range("F30").Value = 38

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 317
$ws.Range("F5").Value = 178
$ws.Range("F8").Value = 391

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9512
$ws.Range("F4").Value = 639
$ws.Range("F5").Value = 210

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9512
$ws.Range("F4").Value = 639
$ws.Range("F5").Value = 372
$ws.Range("F6").Value = 426
$ws.Range("F7").Value = 1155
$ws.Range("F11").Value = 317
$ws.Range("F12").Value = 210
$ws.Range("F13").Value = 1635
$ws.Range("F14").Value = 6130
$ws.Range("G14").Value = 68
$ws.Range("F16").Value = 1770
$ws.Range("F18").Value = 459
$ws.Range("F20").Value = 1
$ws.Range("F22").Value = 6060
$ws.Range("F23").Value = 6060
$ws.Range("F27").Value = 99
$ws.Range("F28").Value = 1668
$ws.Range("F32").Value = 1448
$ws.Range("F33").Value = 733
$ws.Range("F35").Value = 270

